$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Row 7: "Experimental" -> Value cell (B7) was empty, now becomes the literal
#     text "false" (a shared string, NOT a boolean cell).
#     A direct Range.Value assignment of "false"/"FALSE" is auto-coerced by the
#     engine into a genuine Boolean cell (t="b"), which does not match the
#     desired t="s" shared-string cell. Routing the text through a formula
#     ( =T("false") ) and then converting the formula to a static value via
#     Copy / PasteSpecial(xlPasteValues) preserves it as plain text and keeps
#     the existing cell style untouched.
$expCell = $ws.Range("B7")
$expCell.Formula = '=T("false")'
$expCell.Copy()
$expCell.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# --- Row 8: "Date" -> Value cell (B8) timestamp refreshed.
$ws.Range("B8").Value = "2025-06-13T15:45:04+00:00"
